# Weekly data refresh: insert a new price record at the top of the
# "Packham's Triumph" block (row 408), pushing the existing rows
# (408-426) down by one (to 409-427).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 408; existing rows 408:426 shift to 409:427.
$ws.Rows.Item(408).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A408").Value = 10
$ws.Range("B408").Value = "Vega Modelo de Temuco"
$ws.Range("C408").Value = "La Araucanía"
$ws.Range("D408").Value = 44461
$ws.Range("E408").Value = 9
$ws.Range("F408").Value = "Fruta"
$ws.Range("G408").Value = 100104
$ws.Range("H408").Value = "Frutos de pepita"
$ws.Range("I408").Value = 100104005
$ws.Range("J408").Value = "Pera"
$ws.Range("K408").Value = "Packham's Triumph"
$ws.Range("L408").Value = "Primera"
$ws.Range("M408").Value = 160
$ws.Range("N408").Value = 12000
$ws.Range("O408").Value = 13000
$ws.Range("P408").Value = 12500
$ws.Range("Q408").Value = "$/bandeja 18 kilos granel"
$ws.Range("R408").Value = "Región de O'Higgins"
$ws.Range("S408").Value = 694
$ws.Range("T408").Value = 18
